$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a weekly-repeating block of 4 rows (Especial/Primera/
# Segunda/Tercera) for "Pepino dulce" at Lo Valledor de Santiago. A new
# week's block (fecha 44785) needs to be inserted right before the existing
# row 251 block, pushing all the following rows down by 4 (old R299 -> R303).

# Insert 4 blank rows at 251, shifting rows 251:299 down to 255:303.
$ws.Rows("251:254").Insert()

# Shared values for this product across the whole column (unchanged by the edit).
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112043
$categoria = "Pepino dulce"
$variedad = "Cultivar IV Región"
$unidad = "`$/bandeja 18 kilos"
$origen = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"
$fecha = 44785

# New block of 4 rows for the new "fecha" (44785), one per calidad.
$rows = @(
    @{ Row = 251; Calidad = "Especial"; Volumen = 190; PMin = 16000; PMax = 16000; PProm = 16000; PKg = 889 },
    @{ Row = 252; Calidad = "Primera";  Volumen = 330; PMin = 14000; PMax = 14000; PProm = 14000; PKg = 778 },
    @{ Row = 253; Calidad = "Segunda";  Volumen = 240; PMin = 11000; PMax = 11000; PProm = 11000; PKg = 611 },
    @{ Row = 254; Calidad = "Tercera";  Volumen = 210; PMin = 8000;  PMax = 8000;  PProm = 8000;  PKg = 444 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoriaId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
